# Auto-generated edit script: update IFRS financial data values
# for the "company_list" sheet, rows 2-9, columns D-AJ.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 14602
$ws.Range("E2").Value = 957
$ws.Range("F2").Value = 957
$ws.Range("G2").Value = 1036
$ws.Range("H2").Value = 811
$ws.Range("I2").Value = 813
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 14030
$ws.Range("L2").Value = 4200
$ws.Range("M2").Value = 9830
$ws.Range("N2").Value = 9824
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 1462
$ws.Range("Q2").Value = 1025
$ws.Range("R2").Value = -512
$ws.Range("S2").Value = 40
$ws.Range("T2").Value = 347
$ws.Range("U2").Value = 678
$ws.Range("V2").Value = 1978
$ws.Range("W2").Value = 6.55
$ws.Range("X2").Value = 5.55
$ws.Range("Y2").Value = 8.58
$ws.Range("Z2").Value = 6.02
$ws.Range("AA2").Value = 42.73
$ws.Range("AB2").Value = 582.16
$ws.Range("AC2").Value = 2782
$ws.Range("AD2").Value = 10.75
$ws.Range("AE2").Value = 33599
$ws.Range("AF2").Value = 0.89
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.67
$ws.Range("AI2").Value = 17.97
$ws.Range("AJ2").Value = 29240000
$ws.Range("D3").Value = 15710
$ws.Range("E3").Value = 741
$ws.Range("F3").Value = 741
$ws.Range("G3").Value = 734
$ws.Range("H3").Value = 506
$ws.Range("I3").Value = 509
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 14749
$ws.Range("L3").Value = 4573
$ws.Range("M3").Value = 10176
$ws.Range("N3").Value = 10173
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 1462
$ws.Range("Q3").Value = 517
$ws.Range("R3").Value = -1573
$ws.Range("S3").Value = -16
$ws.Range("T3").Value = 337
$ws.Range("U3").Value = 180
$ws.Range("V3").Value = 2306
$ws.Range("W3").Value = 4.72
$ws.Range("X3").Value = 3.22
$ws.Range("Y3").Value = 5.09
$ws.Range("Z3").Value = 3.52
$ws.Range("AA3").Value = 44.94
$ws.Range("AB3").Value = 606.61
$ws.Range("AC3").Value = 1742
$ws.Range("AD3").Value = 15.5
$ws.Range("AE3").Value = 34790
$ws.Range("AF3").Value = 0.78
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.85
$ws.Range("AI3").Value = 28.71
$ws.Range("AJ3").Value = 29240000
$ws.Range("D4").Value = 15293
$ws.Range("E4").Value = 790
$ws.Range("F4").Value = 790
$ws.Range("G4").Value = 742
$ws.Range("H4").Value = 512
$ws.Range("I4").Value = 517
$ws.Range("J4").Value = -5
$ws.Range("K4").Value = 15811
$ws.Range("L4").Value = 5290
$ws.Range("M4").Value = 10522
$ws.Range("N4").Value = 10528
$ws.Range("O4").Value = -6
$ws.Range("P4").Value = 1462
$ws.Range("Q4").Value = 1413
$ws.Range("R4").Value = -2337
$ws.Range("S4").Value = 528
$ws.Range("T4").Value = 286
$ws.Range("U4").Value = 1127
$ws.Range("V4").Value = 2968
$ws.Range("W4").Value = 5.17
$ws.Range("X4").Value = 3.35
$ws.Range("Y4").Value = 4.99
$ws.Range("Z4").Value = 3.35
$ws.Range("AA4").Value = 50.27
$ws.Range("AB4").Value = 631.79
$ws.Range("AC4").Value = 1767
$ws.Range("AD4").Value = 12.2
$ws.Range("AE4").Value = 36006
$ws.Range("AF4").Value = 0.6
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 2.32
$ws.Range("AI4").Value = 28.3
$ws.Range("AJ4").Value = 29240000
$ws.Range("D5").Value = 16021
$ws.Range("E5").Value = 1101
$ws.Range("F5").Value = 1101
$ws.Range("G5").Value = 1019
$ws.Range("H5").Value = 700
$ws.Range("I5").Value = 736
$ws.Range("J5").Value = -37
$ws.Range("K5").Value = 15948
$ws.Range("L5").Value = 4915
$ws.Range("M5").Value = 11033
$ws.Range("N5").Value = 10925
$ws.Range("O5").Value = 108
$ws.Range("P5").Value = 1462
$ws.Range("Q5").Value = 1195
$ws.Range("R5").Value = -287
$ws.Range("S5").Value = -850
$ws.Range("T5").Value = 1055
$ws.Range("U5").Value = 140
$ws.Range("V5").Value = 2269
$ws.Range("W5").Value = 6.88
$ws.Range("X5").Value = 4.37
$ws.Range("Y5").Value = 6.87
$ws.Range("Z5").Value = 4.41
$ws.Range("AA5").Value = 44.55
$ws.Range("AB5").Value = 672.28
$ws.Range("AC5").Value = 2519
$ws.Range("AD5").Value = 12.43
$ws.Range("AE5").Value = 37363
$ws.Range("AF5").Value = 0.84
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 1.6
$ws.Range("AI5").Value = 19.85
$ws.Range("AJ5").Value = 29240000
$ws.Range("D6").Value = 17067
$ws.Range("E6").Value = 1195
$ws.Range("F6").Value = 1195
$ws.Range("G6").Value = 1202
$ws.Range("H6").Value = 813
$ws.Range("I6").Value = 832
$ws.Range("K6").Value = 16518
$ws.Range("L6").Value = 4684
$ws.Range("M6").Value = 11834
$ws.Range("N6").Value = 11733
$ws.Range("P6").Value = 1462
$ws.Range("Q6").Value = 1247
$ws.Range("R6").Value = 1390
$ws.Range("S6").Value = -536
$ws.Range("T6").Value = 491
$ws.Range("U6").Value = 756
$ws.Range("V6").Value = 1858
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 4.77
$ws.Range("Y6").Value = 7.35
$ws.Range("Z6").Value = 5.01
$ws.Range("AA6").Value = 39.58
$ws.Range("AB6").Value = 717.6900000000001
$ws.Range("AC6").Value = 2846
$ws.Range("AD6").Value = 8.73
$ws.Range("AE6").Value = 40126
$ws.Range("AF6").Value = 0.62
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 2.01
$ws.Range("AI6").Value = 17.57
$ws.Range("AJ6").Value = 29240000
$ws.Range("D7").Value = 18454
$ws.Range("E7").Value = 1074
$ws.Range("G7").Value = 1086
$ws.Range("H7").Value = 771
$ws.Range("I7").Value = 800
$ws.Range("K7").Value = 20303
$ws.Range("L7").Value = 6831
$ws.Range("M7").Value = 13472
$ws.Range("N7").Value = 12316
$ws.Range("P7").Value = 1462
$ws.Range("Q7").Value = 777
$ws.Range("R7").Value = -886
$ws.Range("S7").Value = -536
$ws.Range("T7").Value = 254
$ws.Range("U7").Value = 660
$ws.Range("W7").Value = 5.82
$ws.Range("X7").Value = 4.18
$ws.Range("Y7").Value = 6.65
$ws.Range("Z7").Value = 4.19
$ws.Range("AA7").Value = 50.71
$ws.Range("AC7").Value = 2735
$ws.Range("AD7").Value = 5.5
$ws.Range("AE7").Value = 42119
$ws.Range("AF7").Value = 0.36
$ws.Range("AG7").Value = 500
$ws.Range("AH7").Value = 3.32
$ws.Range("AI7").Value = 18.28
$ws.Range("D8").Value = 19190
$ws.Range("E8").Value = 1382
$ws.Range("G8").Value = 1355
$ws.Range("H8").Value = 1032
$ws.Range("I8").Value = 974
$ws.Range("K8").Value = 20982
$ws.Range("L8").Value = 6620
$ws.Range("M8").Value = 14363
$ws.Range("N8").Value = 13146
$ws.Range("P8").Value = 1462
$ws.Range("Q8").Value = 1211
$ws.Range("R8").Value = -359
$ws.Range("S8").Value = -553
$ws.Range("T8").Value = 258
$ws.Range("U8").Value = 939
$ws.Range("W8").Value = 7.2
$ws.Range("X8").Value = 5.38
$ws.Range("Y8").Value = 7.65
$ws.Range("Z8").Value = 5
$ws.Range("AA8").Value = 46.09
$ws.Range("AC8").Value = 3330
$ws.Range("AD8").Value = 4.52
$ws.Range("AE8").Value = 44958
$ws.Range("AF8").Value = 0.33
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 3.32
$ws.Range("AI8").Value = 15.01
$ws.Range("D9").Value = 20020
$ws.Range("E9").Value = 1557
$ws.Range("G9").Value = 1507
$ws.Range("H9").Value = 1152
$ws.Range("I9").Value = 1107
$ws.Range("K9").Value = 21940
$ws.Range("L9").Value = 6567
$ws.Range("M9").Value = 15373
$ws.Range("N9").Value = 14086
$ws.Range("P9").Value = 1462
$ws.Range("Q9").Value = 1291
$ws.Range("R9").Value = -304
$ws.Range("S9").Value = -417
$ws.Range("T9").Value = 99
$ws.Range("U9").Value = 1402
$ws.Range("W9").Value = 7.78
$ws.Range("X9").Value = 5.75
$ws.Range("Y9").Value = 8.130000000000001
$ws.Range("Z9").Value = 5.37
$ws.Range("AA9").Value = 42.72
$ws.Range("AC9").Value = 3785
$ws.Range("AD9").Value = 3.98
$ws.Range("AE9").Value = 48175
$ws.Range("AF9").Value = 0.31
$ws.Range("AG9").Value = 500
$ws.Range("AH9").Value = 3.32
$ws.Range("AI9").Value = 13.21
